Write-Output (Test-Path "C:\Temp\test.txt")
[System.IO.File]::WriteAllText("C:\Temp\test2.txt", "hi")
Write-Output ([System.IO.File]::Exists("C:\Temp\test2.txt"))
Write-Output ([System.IO.File]::ReadAllText("C:\Temp\test2.txt"))
